$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 0.7948164146868251
$ws.Range("F2").Value = 0.8479481641468682
$ws.Range("E3").Value = 0.7697624190064795
$ws.Range("E4").Value = 0.7680345572354211
$ws.Range("F4").Value = 0.824622030237581
$ws.Range("F5").Value = 0.8215982721382289
